$d = $word.ActiveDocument

$replacements = @(
    @("2023-08-06 Sunday", "2023-08-07 Monday"),
    @("11÷9=", "59÷2="),
    @("96÷6=", "59÷2="),
    @("23÷6=", "60÷2="),
    @("19÷4=", "31÷5="),
    @("41÷4=", "62÷2="),
    @("88÷2=", "79÷7="),
    @("39÷9=", "34÷5="),
    @("95÷8=", "54÷9="),
    @("38÷3=", "74÷3="),
    @("42÷6=", "45÷8="),
    @("72÷5=", "99÷6="),
    @("89÷5=", "31÷8="),
    @("11÷4=", "23÷4="),
    @("68÷2=", "22÷6="),
    @("85÷9=", "87÷4="),
    @("39÷3=", "13÷4="),
    @("71÷8=", "69÷9="),
    @("89÷6=", "10÷7="),
    @("39÷5=", "18÷2="),
    @("34÷9=", "66÷3="),
    @("50÷8=", "83÷6="),
    @("65÷2=", "94÷7="),
    @("54÷3=", "74÷3="),
    @("67÷6=", "15÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
